$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.657.12'
$ws.Range('E2').Value = '  +0.04%  '

# Row 3
$ws.Range('D3').Value = '2.294.03'
$ws.Range('E3').Value = '  +0.29%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +18.70%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '268.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.30%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.625'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.43%  '

# Row 8
$ws.Range('E8').Value = '  +0.19%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.624'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.33%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.28'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.40%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0949'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.62%  '

# Row 12
$ws.Range('E12').Value = '  +11.69%  '

# Row 13
$ws.Range('E13').Value = '  +0.90%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.72%  '

# Row 15
$ws.Range('D15').Value = '2.635.70'
$ws.Range('E15').Value = '  +0.21%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.848'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.04%  '

# Row 17
$ws.Range('D17').Value = '2.290.92'
$ws.Range('E17').Value = '  +0.08%  '

# Row 18
$ws.Range('D18').Value = '43.681.82'
$ws.Range('E18').Value = '  +0.20%  '

# Row 19
$ws.Range('E19').Value = '  +2.47%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.71%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.36%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.49'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.80%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.06%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.96%  '

# Row 25
$ws.Range('E25').Value = '  +10.83%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.07%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.08%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.85%  '

# Row 29
$ws.Range('E29').Value = '  -1.87%  '

# Row 30
$ws.Range('E30').Value = '  -0.88%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '177.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.07%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.61'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.02%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0930'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.05%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.61'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.63%  '

# Row 35
$ws.Range('E35').Value = '  +0.93%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.76'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.49%  '

# Row 37
$ws.Range('E37').Value = '  +2.15%  '

# Row 38
$ws.Range('E38').Value = '  +0.63%  '

# Row 39
$ws.Range('E39').Value = '  +10.76%  '

# Row 40
$ws.Range('E40').Value = '  +14.30%  '

# Row 41
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.244'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.16%  '

# Row 42
$ws.Range('B42').Value = 'LidoDAOToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.41'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.38%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.62'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.36%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +18.01%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.30%  '

# Row 46
$ws.Range('E46').Value = '  +0.03%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.20%  '

# Row 48
$ws.Range('E48').Value = '  -0.76%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.70%  '

# Row 50
$ws.Range('E50').Value = '  +3.33%  '

# Row 51
$ws.Range('E51').Value = '  +5.55%  '
